$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "safe_column_name_test" worksheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "safe_column_name_test"

# --- Populate the new sheet's data (8 cols x 6 rows, header + 5 data rows) ---
$data = New-Object 'object[,]' 6,8
$data[0,0] = "col1"
$data[0,1] = "  col with leading and trailing spaces.  "
$data[0,2] = "123_starts_with_123"
$data[0,3] = "With * / special ? ! Characters. "
$data[0,4] = "col1"
$data[0,5] = "col1"
$data[0,6] = "The quick BROWN fox with a very long column name is now jumping over a lazy dog by the zigzag quarry site"
$data[0,7] = "!!!date???"

$data[1,0] = 1
$data[1,1] = "a"
$data[1,2] = "a"
$data[1,3] = 1.5
$data[1,4] = 5
$data[1,5] = "e"
$data[1,6] = "   This is some text. With whitespaces.  "
$data[1,7] = 37145

$data[2,0] = 2
$data[2,1] = "b"
$data[2,2] = "ba"
$data[2,3] = 2.3
$data[2,4] = 4
$data[2,5] = "d"
$data[2,6] = "jumped over the lazy dog"
$data[2,7] = 25023

$data[3,0] = 3
$data[3,1] = "c"
$data[3,2] = "ka"
$data[3,3] = 3.4
$data[3,4] = 3
$data[3,5] = "c"
$data[3,6] = "     by the zigzag`nquarry site.   "
$data[3,7] = "not a date"

$data[4,0] = 4
$data[4,1] = "d"
$data[4,2] = "da"
$data[4,3] = 3.14
$data[4,4] = 2
$data[4,5] = "b"
$data[4,6] = "lorem ipsum dolorem"
$data[4,7] = 1035

$data[5,0] = 5
$data[5,1] = "e"
$data[5,2] = "e"
$data[5,3] = 0.00012
$data[5,4] = 1
$data[5,5] = "a"
$data[5,6] = "Joel was here"
$data[5,7] = 654.34

$ws2.Range("A1:H6").Value2 = $data

# Date-like formatting on the H column ("!!!date???")
$ws2.Range("H2").NumberFormat = "mm-dd-yy"       # -> builtin numFmtId 14 (m/d/yyyy)
$ws2.Range("H3").NumberFormat = "d-mmm-yy"       # -> builtin numFmtId 15 (shared w/ Sheet1 dates)

# Wrapped multi-line text cell
$ws2.Range("G4").WrapText = $true
$ws2.Range("A4").RowHeight = 32

# Column widths to match the source column widths
$ws2.Columns.Item(2).ColumnWidth = 32.5
$ws2.Columns.Item(3).ColumnWidth = 19
$ws2.Columns.Item(4).ColumnWidth = 26.33
$ws2.Columns.Item(7).ColumnWidth = 82.5
$ws2.Columns.Item(8).ColumnWidth = 18.33

# Select G2 and make this newly-added sheet the active/visible tab,
# matching the authored workbook's final UI state.
$ws2.Range("G2").Select() | Out-Null
$ws2.Activate() | Out-Null

# --- Sheet1 cosmetic cleanup: rows 6-9's label cells drop their stray
#     fill/border-applying style in favour of the same plain "Consolas,
#     vertical-center" style already used by rows 2-4 ---
foreach ($r in 6, 7, 8, 9) {
    $cell = $ws1.Range("A$r")
    $cell.Font.Name = "Consolas"
    $cell.Font.Color = 7901646
}
